$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) First bullet: merge "Open your circuit" + " schematic" runs into a
#    single run reading "Move your circuit schematic to this folder".
# -----------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$bodyRange = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$bodyRange.Text = "Move your circuit schematic to this folder"

# -----------------------------------------------------------------------
# 2) Insert a brand-new bullet right after it, reading "Open your" +
#    " schematic" (kept as two separate runs, matching the target XML).
#    We inject raw OOXML via InsertXML so the run split is preserved
#    exactly (plain Range.Text/InsertAfter calls get silently merged
#    into a single run by the editing engine).
# -----------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$insertionPoint = $d.Range($p1.Range.End, $p1.Range.End)

$newParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Open your</w:t></w:r><w:r><w:t xml:space="preserve"> schematic</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$insertionPoint.InsertXML($newParaXml)

# InsertXML leaves behind a stray empty paragraph (used only to force the
# paragraph break) - remove it now that the break has been created.
$strayEmptyPara = $d.Paragraphs(3)
$strayEmptyPara.Range.Delete()

# -----------------------------------------------------------------------
# 3) Resize the first picture (module image) and bump its edit id /
#    effect extent, matching the target drawing XML exactly.
# -----------------------------------------------------------------------
$pic = $d.InlineShapes.Item(1)
$picRange = $pic.Range

$newDrawingXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing"><w:body><w:p><w:r><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="2FFA1FF8" wp14:editId="474E2760"><wp:extent cx="1485995" cy="1593850"/><wp:effectExtent l="0" t="0" r="0" b="6350"/><wp:docPr id="591550916" name="Picture 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 1"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId5"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="1502241" cy="1611275"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$picRange.InsertXML($newDrawingXml)

# The insert above places the resized picture as a *new* run right next
# to the original one; delete the original (still item #1) so only the
# resized picture remains.
$d.InlineShapes.Item(1).Delete()
